# Updates cryptocurrency price/volume figures on the "cryptos" sheet to
# reflect the latest scrape (commit: "Updated symbol list on Fri Jan 20
# 20:28:47 UTC 2023 with GitHub Actions").
#
# All Price (column D) and Volume(1h) (column E) values in this sheet are
# stored as plain text (e.g. "295.25", "0.46%"), not numbers, so each write
# below forces the cell to Text format before assigning the new string and
# then clears the temporary number-format override (ClearFormats) so the
# cell's style index is left exactly as it was (unstyled / General), only
# the stored value changes - matching the source workbook where these rows
# never carried an explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$NewValue)

    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "296.41"
Set-TextValue $ws "E2" "0.73%"
Set-TextValue $ws "D3" "31.67"
Set-TextValue $ws "E3" "2.06%"
Set-TextValue $ws "D4" "4.970"
Set-TextValue $ws "E4" "0.63%"
Set-TextValue $ws "E5" "3.77%"
Set-TextValue $ws "D6" "2.237"
Set-TextValue $ws "E6" "-4.78%"
Set-TextValue $ws "D7" "7.859"
Set-TextValue $ws "E7" "1.78%"
Set-TextValue $ws "D8" "0.9275"
Set-TextValue $ws "E8" "2.76%"
Set-TextValue $ws "D9" "0.09697"
Set-TextValue $ws "E9" "21.47%"
Set-TextValue $ws "D10" "0.1743"
Set-TextValue $ws "E10" "3.53%"
Set-TextValue $ws "D11" "0.08417"
Set-TextValue $ws "E11" "2.80%"
Set-TextValue $ws "D12" "0.03241"
Set-TextValue $ws "E12" "4.47%"
Set-TextValue $ws "D13" "0.09878"
Set-TextValue $ws "E13" "-2.03%"
Set-TextValue $ws "D14" "0.001471"
Set-TextValue $ws "E14" "-1.87%"
Set-TextValue $ws "D15" "0.005683"
Set-TextValue $ws "E15" "-2.24%"
Set-TextValue $ws "D16" "3.495"
Set-TextValue $ws "E16" "0.41%"
Set-TextValue $ws "D17" "3.778"
Set-TextValue $ws "E17" "0.99%"
Set-TextValue $ws "D18" "2.196"
Set-TextValue $ws "E18" "5.89%"
Set-TextValue $ws "D19" "0.3353"
Set-TextValue $ws "E19" "0.64%"
Set-TextValue $ws "D20" "0.1321"
Set-TextValue $ws "E20" "1.45%"
Set-TextValue $ws "D21" "4.074"
Set-TextValue $ws "E21" "2.12%"
Set-TextValue $ws "D22" "0.2282"
Set-TextValue $ws "E22" "8.85%"
Set-TextValue $ws "D23" "0.04520"
Set-TextValue $ws "E23" "-0.17%"
Set-TextValue $ws "D24" "0.001209"
Set-TextValue $ws "E24" "-0.10%"
Set-TextValue $ws "D25" "0.004352"
Set-TextValue $ws "E25" "-6.57%"
Set-TextValue $ws "D26" "0.0001291"
Set-TextValue $ws "E26" "-0.47%"
Set-TextValue $ws "D27" "0.0003362"
Set-TextValue $ws "E27" "-0.75%"
Set-TextValue $ws "D39" "0.01679"
Set-TextValue $ws "E39" "4.34%"
Set-TextValue $ws "D40" "0.04627"
Set-TextValue $ws "E40" "3.92%"
Set-TextValue $ws "D41" "0.007512"
Set-TextValue $ws "E41" "2.04%"
Set-TextValue $ws "D42" "0.009739"
Set-TextValue $ws "E42" "13.23%"
Set-TextValue $ws "D43" "0.1387"
Set-TextValue $ws "E43" "4.50%"
Set-TextValue $ws "D44" "0.002144"
Set-TextValue $ws "E44" "7.31%"
Set-TextValue $ws "E45" "-0.77%"
Set-TextValue $ws "D46" "0.00006051"
Set-TextValue $ws "E46" "2.50%"
Set-TextValue $ws "D47" "0.00000000745"
Set-TextValue $ws "E47" "-0.45%"
Set-TextValue $ws "D48" "2.794"
Set-TextValue $ws "E48" "24.69%"
Set-TextValue $ws "D49" "0.001981"
Set-TextValue $ws "E49" "-31.51%"
Set-TextValue $ws "D50" "0.00002086"
Set-TextValue $ws "E50" "-0.45%"
Set-TextValue $ws "D51" "0.0001987"
Set-TextValue $ws "E51" "-0.45%"
